# Update the "2^40,000,000,000" Fibonacci-power record on slide 1 to
# "2^50,000,000,000" without disturbing any other run/paragraph in the
# shape (the slide also contains an unrelated "40,000,000,000" — the
# "...the 40,000,000,000th Fibonacci number" line — which must stay
# untouched).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the body/content placeholder shape. Prefer matching by name,
# but fall back to scanning every shape with a text frame in case the
# name ever differs.
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Content Placeholder 2") {
        $shape = $candidate
    }
}

if ($null -eq $shape) {
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $candidate = $s.Shapes.Item($i)
        if ($candidate.HasTextFrame -and $candidate.TextFrame.HasText) {
            if ($candidate.TextFrame.TextRange.Text -like "*40,000,000,000*") {
                $shape = $candidate
            }
        }
    }
}

if ($null -eq $shape) {
    throw "Could not find the shape containing the Fibonacci-power record."
}

$tr = $shape.TextFrame.TextRange

# The old/new values, plus enough surrounding text to uniquely identify
# the "Computed  2^40,000,000,000." run (as opposed to the other,
# unrelated "40,000,000,000" occurring earlier in the same shape).
$oldNumber = "40,000,000,000"
$newNumber = "50,000,000,000"
$needle = "Computed  2" + $oldNumber + "."

$needleLen = $needle.Length
$totalLen = $tr.Length

# `TextRange.Characters(start, length)` is used instead of whole-range
# `.Text` assignment so the existing run/paragraph structure (fonts,
# the superscript "th" run, etc.) is preserved; only the digits change.
$foundAt = -1
for ($i = 1; $i -le ($totalLen - $needleLen + 1); $i++) {
    $chunk = $tr.Characters($i, $needleLen).Text
    if ($chunk -eq $needle) {
        $foundAt = $i
        break
    }
}

if ($foundAt -lt 0) {
    throw "Could not locate the 'Computed  2^40,000,000,000.' run to update."
}

$numberStart = $foundAt + "Computed  2".Length
$target = $tr.Characters($numberStart, $oldNumber.Length)
$target.Text = $newNumber
